# Weekly update: insert a new price record right after the header's first
# data row (pushing all existing records down by one) and append a new
# record at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert a new row at row 3 (shifts old rows 3..22 down to 4..23) ---
$ws.Rows.Item(3).Insert()

# --- Step 2: populate the newly inserted row 3 with this week's record ---
$ws.Range("A3").Value = 7
$ws.Range("B3").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C3").Value = "Ñuble"
$ws.Range("D3").Value = 44831
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 300000000
$ws.Range("G3").Value = "Espárragos"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 2200
$ws.Range("M3").Value = 2100
$ws.Range("N3").Value = "$/kilo"
$ws.Range("O3").Value = "Provincia de Diguillín"
$ws.Range("P3").Value = 2100
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = "Hortaliza"

# --- Step 3: append a new record as row 24 (after the shift, data now ends at row 23) ---
$ws.Range("A24").Value = 7
$ws.Range("B24").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C24").Value = "Ñuble"
$ws.Range("D24").Value = 44832
$ws.Range("E24").Value = 16
$ws.Range("F24").Value = 300000000
$ws.Range("G24").Value = "Espárragos"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 200
$ws.Range("K24").Value = 2000
$ws.Range("L24").Value = 2200
$ws.Range("M24").Value = 2100
$ws.Range("N24").Value = "$/kilo"
$ws.Range("O24").Value = "Provincia de Diguillín"
$ws.Range("P24").Value = 2100
$ws.Range("Q24").Value = 1
$ws.Range("R24").Value = "Hortaliza"

# Match the date-formatted number format used by the rest of column D
$ws.Range("D24").NumberFormat = "YYYY-MM-DD HH:MM:SS"
